$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Update system quantity driver B1: 2 -> 25 (all dependent formulas recalc automatically)
$ws.Range("B1").Value = 25

# Fill in the part number for the red LED (D1) row - was previously blank
$ws.Range("B20").Value = "TLMS1000-GS08"

# Add a new BOM line for the male 100 mil headers (J2), flagged in red since it
# hasn't been priced/ordered yet
$ws.Range("A28").Value = "100 mil headers, male"
$ws.Range("B28").Value = "J2"
$ws.Range("A28:B28").Interior.Color = 255

# Re-fit the Name/Designator columns now that B20 has new (wider) content
$ws.Columns.Item(2).ColumnWidth = 20.2
$ws.Columns.Item(3).ColumnWidth = 13.3

# Leave the cursor where the new part number was typed in
$ws.Range("B16").Select()

Write-Output "done"
